# Oracle India jobs sheet update:
#  - Two brand new "Human Resources Associate" postings are inserted at the
#    top of the data (new rows 2 and 3), pushing the existing three rows
#    down by two positions (old rows 2-4 become rows 4-6, unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 2, shifting existing data down.
$ws.Range("A2:E3").Insert()

# The inserted cells pick up formatting from the row above (the bold header);
# strip that back to the default/unstyled look used by the other data rows.
$ws.Range("A2:E3").ClearFormats()

# --- New row 2: Human Resources Associate / job 31613 ---
$ws.Range("A2").Value = "Oracle Careers"
$ws.Range("B2").Value = "Human Resources Associate"
$ws.Range("C2").Value = "New Delhi, India"
# Force the posting date to stay literal text (as in the source data)
# instead of being auto-parsed into a date serial number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "02/01/2026"
$ws.Range("E2").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/31613/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# --- New row 3: Human Resources Associate / job 31614 ---
$ws.Range("A3").Value = "Oracle Careers"
$ws.Range("B3").Value = "Human Resources Associate"
$ws.Range("C3").Value = "New Delhi, India"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "02/01/2026"
$ws.Range("E3").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/31614/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# The date-formatted NumberFormat picked up by the inserted cells needs to be
# stripped back to the unstyled default once the literal text values (not
# real dates) are committed.
$ws.Range("D2:D3").ClearFormats()
